$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Report/Print date+time text (row 7, static text cells) ---
# Stored with a leading apostrophe in the original file (quotePrefix) so the
# dd/mm/yyyy-looking text isn't reinterpreted as a real date/time value.
$ws.Range("G7").Value = "'20/04/2020"
$ws.Range("H7").Value = "'17:19:00"

# --- Designator column (J11:J26) renumbering ---
# Leading apostrophe preserves the original quotePrefix-flagged cell style.
$ws.Range("J11").Value = "'R2, R3, R5, R6, R7, R8, R9, R10, R11, R12, R13, R14, R15"
$ws.Range("J15").Value = "'P3, P4"
$ws.Range("J17").Value = "'U1"
$ws.Range("J18").Value = "'R4"
$ws.Range("J19").Value = "'C1"
$ws.Range("J20").Value = "'C2"
$ws.Range("J21").Value = "'C3"
$ws.Range("J23").Value = "'P2"
$ws.Range("J24").Value = "'P1"
$ws.Range("J25").Value = "'U2"
